$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "0.9991", "42.08") are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.490.11'
$ws.Range("E2").Value = '  -3.22%  '

# Row 3
$ws.Range("D3").Value = '1.752.20'
$ws.Range("E3").Value = '  -3.68%  '

# Row 4
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").Value = '322.23'
$ws.Range("E5").Value = '  -2.25%  '

# Row 6
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").Value = '0.4250'
$ws.Range("E7").Value = '  -4.60%  '

# Row 8
$ws.Range("D8").Value = '0.3606'
$ws.Range("E8").Value = '  -3.15%  '

# Row 9
$ws.Range("D9").Value = '0.07482'
$ws.Range("E9").Value = '  -3.46%  '

# Row 10
$ws.Range("D10").Value = '42.08'
$ws.Range("E10").Value = '  -6.34%  '

# Row 11
$ws.Range("D11").Value = '1.101'
$ws.Range("E11").Value = '  -2.92%  '

# Row 12
$ws.Range("D12").Value = '0.9976'
$ws.Range("E12").Value = '  -0.14%  '

# Row 13
$ws.Range("D13").Value = '20.72'
$ws.Range("E13").Value = '  -6.70%  '

# Row 14
$ws.Range("D14").Value = '6.034'
$ws.Range("E14").Value = '  -4.51%  '

# Row 15
$ws.Range("D15").Value = '7.237'
$ws.Range("E15").Value = '  -4.92%  '

# Row 16
$ws.Range("D16").Value = '1.744.50'
$ws.Range("E16").Value = '  -5.99%  '

# Row 17
$ws.Range("D17").Value = '92.73'
$ws.Range("E17").Value = '  -0.86%  '

# Row 18
$ws.Range("D18").Value = '0.00001065'
$ws.Range("E18").Value = '  -1.93%  '

# Row 19
$ws.Range("D19").Value = '0.06377'
$ws.Range("E19").Value = '  -2.41%  '

# Row 20
$ws.Range("D20").Value = '0.9986'
$ws.Range("E20").Value = '  -0.03%  '

# Row 21
$ws.Range("D21").Value = '17.06'
$ws.Range("E21").Value = '  -3.01%  '

# Row 22
$ws.Range("D22").Value = '5.891'
$ws.Range("E22").Value = '  -5.87%  '

# Row 23
$ws.Range("D23").Value = '27.509.24'
$ws.Range("E23").Value = '  -3.26%  '

# Row 24
$ws.Range("D24").Value = '11.24'
$ws.Range("E24").Value = '  -4.06%  '

# Row 25
$ws.Range("D25").Value = '2.101'
$ws.Range("E25").Value = '  -3.97%  '

# Row 26
$ws.Range("D26").Value = '160.88'
$ws.Range("E26").Value = '  +3.06%  '

# Row 27
$ws.Range("D27").Value = '20.29'
$ws.Range("E27").Value = '  -2.68%  '

# Row 28
$ws.Range("D28").Value = '1.943.87'
$ws.Range("E28").Value = '  -5.04%  '

# Row 29
$ws.Range("D29").Value = '2.131'
$ws.Range("E29").Value = '  -8.11%  '

# Row 30
$ws.Range("D30").Value = '123.85'
$ws.Range("E30").Value = '  -3.74%  '

# Row 31
$ws.Range("D31").Value = '1.105'
$ws.Range("E31").Value = '  -8.82%  '

# Row 32
$ws.Range("D32").Value = '3.654'
$ws.Range("E32").Value = '  -0.08%  '

# Row 33
$ws.Range("D33").Value = '5.552'
$ws.Range("E33").Value = '  -6.50%  '

# Row 34
$ws.Range("D34").Value = '0.08864'
$ws.Range("E34").Value = '  -4.25%  '

# Row 35
$ws.Range("D35").Value = '12.23'
$ws.Range("E35").Value = '  -7.46%  '

# Row 36
$ws.Range("D36").Value = '0.02292'
$ws.Range("E36").Value = '  -2.95%  '

# Row 37
$ws.Range("E37").Value = '  -4.30%  '

# Row 38
$ws.Range("D38").Value = '0.06000'
$ws.Range("E38").Value = '  -4.01%  '

# Row 39
$ws.Range("D39").Value = '0.6329'
$ws.Range("E39").Value = '  -4.47%  '

# Row 40
$ws.Range("D40").Value = '4.941'
$ws.Range("E40").Value = '  -5.23%  '

# Row 41
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -1.49%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '7.919'
$ws.Range("E42").Value = '  -3.12%  '

# Row 43
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '0.9978'
$ws.Range("E43").Value = '  -0.05%  '

# Row 44
$ws.Range("E44").Value = '  -2.19%  '

# Row 45
$ws.Range("D45").Value = '13.31'
$ws.Range("E45").Value = '  -5.09%  '

# Row 46
$ws.Range("D46").Value = '0.5887'
$ws.Range("E46").Value = '  -4.14%  '

# Row 47
$ws.Range("D47").Value = '3.691'

# Row 48
$ws.Range("D48").Value = '123.25'
$ws.Range("E48").Value = '  -3.34%  '

# Row 49
$ws.Range("D49").Value = '1.969'
$ws.Range("E49").Value = '  -3.82%  '

# Row 50
$ws.Range("D50").Value = '1.166'
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$ws.Range("D51").Value = '0.06825'
$ws.Range("E51").Value = '  -2.55%  '

# Restore default (Normal) style on the price column so no stray
# number-format override is left on cells (matches original formatting).
$priceRange.Style = "Normal"
